$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961433075938"
$ws1.Range("B2").Value = "go_stims-16509961432755675.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996143291593.csv"
$ws1.Range("B4").Value = "go_stims-1650996143291593.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961433075938.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961467395601"
$ws2.Range("B2").Value = "OB-16509961457715664.csv"
$ws2.Range("B3").Value = "OB-16509961445155635.csv"
$ws2.Range("B4").Value = "ZB-match_4-16509961442835598.csv"
$ws2.Range("B5").Value = "ZB-match_6-16509961433875997.csv"
$ws2.Range("B6").Value = "OB-1650996145619558.csv"
$ws2.Range("B7").Value = "TB-16509961465476.csv"
$ws2.Range("B8").Value = "ZB-match_5-16509961437875583.csv"
$ws2.Range("B9").Value = "TB-16509961459715836.csv"
$ws2.Range("B10").Value = "TB-1650996146723594.csv"

# --- Sheet 3: RS_TO (name only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961467395601"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961467875931"
$ws4.Range("B2").Value = "MM_stims-16509961467555654.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961467395601.csv"
$ws4.Range("B4").Value = "MM_stims-16509961467715642.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961467555654.csv"
$ws4.Range("B6").Value = "MM_stims-16509961467875931.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961467715642.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650996146867602"
$ws5.Range("B2").Value = "SAT_stims-16509961468196251.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961468355966.csv"
$ws5.Range("B4").Value = "vSAT_stims-165099614685156.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961467875931.csv"
